$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2875388.8
$ws.Range("I17").Value = 300
$ws.Range("J17").Value = 2978070.5
$ws.Range("K17").Value = 900
$ws.Range("L17").Value = 8934211.5
$ws.Range("M17").Value = -732
$ws.Range("N17").Value = -8934547.5

$ws.Range("H112").Value = 2292.25
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 2292.25
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 6876.75
$ws.Range("N112").Value = -9092.75
$ws.Range("M112").ClearContents()

$ws.Range("H137").Value = 1470.32
$ws.Range("I137").Value = 774.1177
$ws.Range("J137").Value = 2949.75
$ws.Range("K137").Value = 2322.3531
$ws.Range("L137").Value = 8849.25
$ws.Range("M137").Value = 227.6468999999997
$ws.Range("N137").Value = -13949.25

$ws.Range("H138").Value = 1906.6
$ws.Range("I138").Value = 1113.0571
$ws.Range("J138").Value = 3295.3
$ws.Range("K138").Value = 3339.1713
$ws.Range("L138").Value = 9885.900000000001
$ws.Range("M138").Value = 1800.8287
$ws.Range("N138").Value = -20165.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1431.5
$ws.Range("I45").Value = 1257.2222
$ws.Range("J45").Value = 3000
$ws.Range("K45").Value = 1257.2222
$ws.Range("L45").Value = 3000
$ws.Range("M45").Value = -880.2221999999999
$ws.Range("N45").Value = -3754

$ws.Range("H61").Value = 1857.2
$ws.Range("I61").Value = 1988.8928
$ws.Range("J61").Value = 1330.4286
$ws.Range("K61").Value = 1988.8928
$ws.Range("L61").Value = 1330.4286
$ws.Range("M61").Value = -1776.8928
$ws.Range("N61").Value = -1754.4286

$ws.Range("H102").Value = 2194.75
$ws.Range("I102").Value = 2194.75
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2194.75
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -572.75
$ws.Range("N102").ClearContents()

$ws.Range("H122").Value = 7465.2
$ws.Range("I122").Value = 7188.2607
$ws.Range("J122").Value = 10650
$ws.Range("K122").Value = 21564.7821
$ws.Range("L122").Value = 31950
$ws.Range("M122").Value = -19114.7821

$ws.Range("H136").Value = 1857.2
$ws.Range("I136").Value = 1988.8928
$ws.Range("J136").Value = 1330.4286
$ws.Range("K136").Value = 5966.678400000001
$ws.Range("L136").Value = 3991.2858
$ws.Range("M136").Value = -3416.678400000001
$ws.Range("N136").Value = -9091.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3232.9092
$ws.Range("I105").Value = 1554.6316
$ws.Range("J105").Value = 4508.4
$ws.Range("K105").Value = 1554.6316
$ws.Range("L105").Value = 4508.4
$ws.Range("M105").Value = 192.3684000000001
$ws.Range("N105").Value = -8002.4

$ws.Range("H132").Value = 58400
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 58400
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 58400
$ws.Range("N132").Value = -68520

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4431.9546
$ws.Range("I31").Value = 5851.357
$ws.Range("J31").Value = 1948
$ws.Range("K31").Value = 5851.357
$ws.Range("L31").Value = 1948
$ws.Range("M31").Value = -5556.357
$ws.Range("N31").Value = -2538

$ws.Range("H34").Value = 4431.9546
$ws.Range("I34").Value = 5851.357
$ws.Range("J34").Value = 1948
$ws.Range("K34").Value = 5851.357
$ws.Range("L34").Value = 1948
$ws.Range("M34").Value = -5649.357
$ws.Range("N34").Value = -2352

$ws.Range("H58").Value = 3226.4
$ws.Range("I58").Value = 1110.25
$ws.Range("J58").Value = 5644.857
$ws.Range("K58").Value = 1110.25
$ws.Range("L58").Value = 5644.857
$ws.Range("M58").Value = -907.25
$ws.Range("N58").Value = -6050.857

$ws.Range("H132").Value = 2739.2942
$ws.Range("I132").Value = 2764.182
$ws.Range("J132").Value = 2693.6667
$ws.Range("K132").Value = 8292.545999999998
$ws.Range("L132").Value = 8081.000100000001
$ws.Range("M132").Value = -5762.545999999998
$ws.Range("N132").Value = -13141.0001

$ws.Range("H134").Value = 918568.25
$ws.Range("I134").Value = 2939.8096
$ws.Range("J134").Value = 4764207.5
$ws.Range("K134").Value = 8819.4288
$ws.Range("L134").Value = 14292622.5
$ws.Range("M134").Value = -6284.4288
$ws.Range("N134").Value = -14297692.5

$ws.Range("H136").Value = 3226.4
$ws.Range("I136").Value = 1110.25
$ws.Range("J136").Value = 5644.857
$ws.Range("K136").Value = 3330.75
$ws.Range("L136").Value = 16934.571
$ws.Range("M136").Value = -780.75
$ws.Range("N136").Value = -22034.571

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 376.05
$ws.Range("I107").Value = 229.28572
$ws.Range("J107").Value = 455.07693
$ws.Range("K107").Value = 687.85716
$ws.Range("L107").Value = 1365.23079
$ws.Range("M107").Value = 1232.14284
$ws.Range("N107").Value = -5205.23079

$ws.Range("H126").Value = 4486
$ws.Range("I126").Value = 215
$ws.Range("J126").Value = 7333.3335
$ws.Range("K126").Value = 645
$ws.Range("L126").Value = 22000.0005
$ws.Range("M126").Value = 4295
$ws.Range("N126").Value = -31880.0005

$ws.Range("H131").Value = 821.09
$ws.Range("I131").Value = 377.9
$ws.Range("J131").Value = 870.3333
$ws.Range("K131").Value = 1133.7
$ws.Range("L131").Value = 2610.9999
$ws.Range("M131").Value = 3906.3
$ws.Range("N131").Value = -12690.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3965.2354
$ws.Range("I102").Value = 4597.643
$ws.Range("J102").Value = 1014
$ws.Range("K102").Value = 4597.643
$ws.Range("L102").Value = 1014
$ws.Range("M102").Value = -2975.643

$ws.Range("H132").Value = 3623.4075
$ws.Range("I132").Value = 3296.7
$ws.Range("J132").Value = 4556.857
$ws.Range("K132").Value = 9890.099999999999
$ws.Range("L132").Value = 13670.571
$ws.Range("M132").Value = -7360.099999999999
$ws.Range("N132").Value = -18730.571

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 995.7727
$ws.Range("I46").Value = 715
$ws.Range("J46").Value = 1058.1666
$ws.Range("K46").Value = 715
$ws.Range("L46").Value = 1058.1666
$ws.Range("M46").Value = -527
$ws.Range("N46").Value = -1434.1666

$ws.Range("H82").Value = 1769.8
$ws.Range("I82").Value = 1263
$ws.Range("J82").Value = 2530
$ws.Range("K82").Value = 1263
$ws.Range("L82").Value = 2530
$ws.Range("M82").Value = -902
$ws.Range("N82").Value = -3252

$ws.Range("H85").Value = 1769.8
$ws.Range("I85").Value = 1263
$ws.Range("J85").Value = 2530
$ws.Range("K85").Value = 1263
$ws.Range("L85").Value = 2530
$ws.Range("M85").Value = -15
$ws.Range("N85").Value = -5026

$ws.Range("H100").Value = 3200
$ws.Range("I100").Value = 3500
$ws.Range("J100").Value = 2600
$ws.Range("K100").Value = 3500
$ws.Range("L100").Value = 2600
$ws.Range("M100").Value = -2959
$ws.Range("N100").Value = -3682

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3031.75
$ws.Range("I96").Value = 1886.1428
$ws.Range("J96").Value = 3503.4707
$ws.Range("K96").Value = 1886.1428
$ws.Range("L96").Value = 3503.4707
$ws.Range("M96").Value = -513.1428000000001
$ws.Range("N96").Value = -6249.4707
